# Auto-generated edits applying the crypto price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "70.700.13"
$ws.Range("E2").Value = "  +0.18%  "
Set-TextCell "D3" "3.518.18"
$ws.Range("E3").Value = "  -1.33%  "
Set-TextCell "D4" "0.999"
$ws.Range("E4").Value = "  -0.11%  "
Set-TextCell "D5" "622.33"
$ws.Range("E5").Value = "  +3.88%  "
Set-TextCell "D6" "171.87"
$ws.Range("E6").Value = "  -0.37%  "
Set-TextCell "D7" "3.511.84"
$ws.Range("E7").Value = "  -1.36%  "
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  +1.23%  "
Set-TextCell "D11" "7.17"
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("E12").Value = "  -0.56%  "
Set-TextCell "D13" "46.20"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("E14").Value = "  -0.36%  "
Set-TextCell "D15" "4.079.59"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("E16").Value = "  +0.54%  "
Set-TextCell "D17" "606.14"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell "D18" "70.811.04"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D19" "3.507.92"
$ws.Range("E19").Value = "  -2.15%  "
$ws.Range("E20").Value = "  +1.71%  "
Set-TextCell "D21" "17.69"
$ws.Range("E21").Value = "  +1.74%  "
Set-TextCell "D22" "0.880"
$ws.Range("E22").Value = "  -0.28%  "
Set-TextCell "D23" "9.09"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D24" "15.55"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D25" "97.26"
$ws.Range("E25").Value = "  +0.41%  "
Set-TextCell "D26" "3.71"
$ws.Range("E26").Value = "  -1.39%  "
Set-TextCell "D27" "1.00"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  -2.33%  "
Set-TextCell "D29" "33.53"
$ws.Range("E29").Value = "  -0.91%  "
Set-TextCell "D30" "9.03"
$ws.Range("E30").Value = "  -0.94%  "
Set-TextCell "D31" "3.00"
$ws.Range("E31").Value = "  -1.62%  "
Set-TextCell "D32" "8.10"
$ws.Range("E33").Value = "  -0.24%  "
Set-TextCell "D34" "6.81"
$ws.Range("E34").Value = "  -4.69%  "
Set-TextCell "D35" "620.90"
$ws.Range("E35").Value = "  -5.98%  "
$ws.Range("B36").Value = "Cosmos"
$ws.Range("C36").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D36" "10.85"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D37" "0.0490"
$ws.Range("E37").Value = "  +2.65%  "
Set-TextCell "D38" "0.0994"
$ws.Range("E38").Value = "  -1.42%  "
Set-TextCell "D39" "56.67"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("E40").Value = "  +0.32%  "
Set-TextCell "D41" "3.40"
$ws.Range("E41").Value = "  -6.73%  "
Set-TextCell "D42" "0.143"
$ws.Range("E42").Value = "  +1.08%  "
Set-TextCell "D43" "3.334.21"
$ws.Range("E43").Value = "  -1.51%  "
Set-TextCell "D44" "0.0₃0724"
$ws.Range("E44").Value = "  +1.76%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell "D45" "0.310"
$ws.Range("E45").Value = "  -3.38%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextCell "D46" "2.92"
$ws.Range("E46").Value = "  -0.60%  "
Set-TextCell "D47" "31.82"
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("E48").Value = "  -5.50%  "
Set-TextCell "D50" "133.89"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D51" "0.155"
$ws.Range("E51").Value = "  +6.41%  "
